$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44313
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 1300
$ws.Range("M2").Value = 1400
$ws.Range("P2").Value = 467
$ws.Range("D3").Value = 44313
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 900
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = 950
$ws.Range("P3").Value = 317
$ws.Range("D4").Value = 44406
$ws.Range("D5").Value = 44383
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 350
$ws.Range("K5").Value = 2800
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = 2886
$ws.Range("P5").Value = 962
$ws.Range("D6").Value = 44277
$ws.Range("J6").Value = 250
$ws.Range("D7").Value = 44341
$ws.Range("K7").Value = 1400
$ws.Range("M7").Value = 1450
$ws.Range("P7").Value = 483
$ws.Range("D8").Value = 44460
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 950
$ws.Range("M8").Value = 975
$ws.Range("P8").Value = 325
$ws.Range("D9").Value = 44299
$ws.Range("J9").Value = 300
$ws.Range("D10").Value = 44299
$ws.Range("I10").Value = "Segunda"
$ws.Range("J10").Value = 250
$ws.Range("K10").Value = 1000
$ws.Range("L10").Value = 1200
$ws.Range("M10").Value = 1100
$ws.Range("P10").Value = 367
$ws.Range("D11").Value = 44217
$ws.Range("J11").Value = 250
$ws.Range("K11").Value = 2900
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = 2950
$ws.Range("P11").Value = 983
$ws.Range("D12").Value = 44376
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 280
$ws.Range("K12").Value = 2400
$ws.Range("L12").Value = 2500
$ws.Range("M12").Value = 2436
$ws.Range("P12").Value = 812
$ws.Range("D13").Value = 44257
$ws.Range("J13").Value = 1500
$ws.Range("K13").Value = 2800
$ws.Range("L13").Value = 3000
$ws.Range("M13").Value = 2900
$ws.Range("P13").Value = 967
$ws.Range("D14").Value = 44327
$ws.Range("I14").Value = "Primera"
$ws.Range("K14").Value = 1400
$ws.Range("L14").Value = 1500
$ws.Range("M14").Value = 1450
$ws.Range("P14").Value = 483
$ws.Range("D15").Value = 44327
$ws.Range("I15").Value = "Segunda"
$ws.Range("K15").Value = 1000
$ws.Range("L15").Value = 1200
$ws.Range("M15").Value = 1100
$ws.Range("P15").Value = 367
$ws.Range("D16").Value = 44175
$ws.Range("J16").Value = 250
$ws.Range("K16").Value = 1800
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = 1900
$ws.Range("P16").Value = 633
$ws.Range("D17").Value = 44364
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 270
$ws.Range("K17").Value = 3400
$ws.Range("L17").Value = 3500
$ws.Range("M17").Value = 3450
$ws.Range("P17").Value = 1150
$ws.Range("D18").Value = 44418
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 2400
$ws.Range("L18").Value = 2500
$ws.Range("M18").Value = 2450
$ws.Range("O18").Value = "Región de Arica y Parinacota"
$ws.Range("P18").Value = 817
$ws.Range("D19").Value = 44412
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 300
$ws.Range("M19").Value = 2900
$ws.Range("P19").Value = 967
$ws.Range("D20").Value = 44308
$ws.Range("J20").Value = 270
$ws.Range("D21").Value = 44322
$ws.Range("J21").Value = 250
$ws.Range("K21").Value = 1400
$ws.Range("L21").Value = 1500
$ws.Range("M21").Value = 1450
$ws.Range("P21").Value = 483
$ws.Range("D22").Value = 44343
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 1500
$ws.Range("M22").Value = 1500
$ws.Range("P22").Value = 500
$ws.Range("D23").Value = 44343
$ws.Range("I23").Value = "Segunda"
$ws.Range("J23").Value = 150
$ws.Range("L23").Value = 1400
$ws.Range("M23").Value = 1400
$ws.Range("P23").Value = 467
$ws.Range("D24").Value = 44356
$ws.Range("I24").Value = "Primera"
$ws.Range("K24").Value = 2400
$ws.Range("L24").Value = 2500
$ws.Range("M24").Value = 2450
$ws.Range("P24").Value = 817
$ws.Range("D25").Value = 44356
$ws.Range("I25").Value = "Segunda"
$ws.Range("J25").Value = 200
$ws.Range("K25").Value = 1800
$ws.Range("L25").Value = 2000
$ws.Range("M25").Value = 1900
$ws.Range("P25").Value = 633
$ws.Range("D26").Value = 44467
$ws.Range("J26").Value = 250
$ws.Range("K26").Value = 800
$ws.Range("L26").Value = 900
$ws.Range("M26").Value = 850
$ws.Range("P26").Value = 283
$ws.Range("D27").Value = 44273
$ws.Range("J27").Value = 250
$ws.Range("K27").Value = 3800
$ws.Range("L27").Value = 4000
$ws.Range("M27").Value = 3900
$ws.Range("P27").Value = 1300
$ws.Range("D28").Value = 44168
$ws.Range("K28").Value = 1800
$ws.Range("L28").Value = 2000
$ws.Range("M28").Value = 1900
$ws.Range("P28").Value = 633
$ws.Range("D29").Value = 44292
$ws.Range("J29").Value = 270
$ws.Range("K29").Value = 2400
$ws.Range("L29").Value = 2500
$ws.Range("M29").Value = 2450
$ws.Range("P29").Value = 817
$ws.Range("D30").Value = 44335
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 250
$ws.Range("L30").Value = 1500
$ws.Range("M30").Value = 1450
$ws.Range("P30").Value = 483
$ws.Range("D31").Value = 44320
$ws.Range("J31").Value = 200
$ws.Range("D32").Value = 44320
$ws.Range("I32").Value = "Segunda"
$ws.Range("J32").Value = 200
$ws.Range("K32").Value = 1000
$ws.Range("L32").Value = 1200
$ws.Range("M32").Value = 1100
$ws.Range("P32").Value = 367
$ws.Range("D33").Value = 44474
$ws.Range("J33").Value = 270
$ws.Range("K33").Value = 1000
$ws.Range("L33").Value = 1200
$ws.Range("M33").Value = 1100
$ws.Range("P33").Value = 367
$ws.Range("D34").Value = 44300
$ws.Range("J34").Value = 160
$ws.Range("K34").Value = 1000
$ws.Range("L34").Value = 1200
$ws.Range("M34").Value = 1100
$ws.Range("O34").Value = "Región de Coquimbo"
$ws.Range("P34").Value = 367
$ws.Range("D35").Value = 44350
$ws.Range("J35").Value = 300
$ws.Range("K35").Value = 1800
$ws.Range("L35").Value = 2000
$ws.Range("M35").Value = 1900
$ws.Range("P35").Value = 633
$ws.Range("D36").Value = 44448
$ws.Range("K36").Value = 1400
$ws.Range("L36").Value = 1500
$ws.Range("M36").Value = 1450
$ws.Range("P36").Value = 483
$ws.Range("D37").Value = 44448
$ws.Range("K37").Value = 1000
$ws.Range("L37").Value = 1200
$ws.Range("M37").Value = 1100
$ws.Range("P37").Value = 367
$ws.Range("D38").Value = 44435
$ws.Range("J38").Value = 270
$ws.Range("D39").Value = 44392
$ws.Range("J39").Value = 200
$ws.Range("K39").Value = 3800
$ws.Range("L39").Value = 4000
$ws.Range("M39").Value = 3900
$ws.Range("P39").Value = 1300
$ws.Range("D40").Value = 44392
$ws.Range("I40").Value = "Segunda"
$ws.Range("J40").Value = 200
$ws.Range("K40").Value = 3200
$ws.Range("L40").Value = 3500
$ws.Range("M40").Value = 3350
$ws.Range("P40").Value = 1117
$ws.Range("D41").Value = 44432
$ws.Range("K41").Value = 1800
$ws.Range("L41").Value = 2000
$ws.Range("M41").Value = 1900
$ws.Range("P41").Value = 633
